# Updated symbol list on Sun Jan  1 23:20:05 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new literal text value. Values that look numeric/percent
# must be forced to Text format before assignment so Excel keeps them as
# literal strings (matching the original inline-string cells) instead of
# re-interpreting them as numbers/percentages.
$updates = [ordered]@{
    'D2' = '244.68'
    'E2' = '-0.80%'
    'D3' = '27.56'
    'E3' = '6.09%'
    'D4' = '5.121'
    'E4' = '0.52%'
    'D5' = '0.05679'
    'E5' = '1.64%'
    'D6' = '6.493'
    'E6' = '0.16%'
    'E7' = '0.66%'
    'D8' = '0.8531'
    'E8' = '0.98%'
    'E9' = '0.21%'
    'D10' = '0.02871'
    'E10' = '1.87%'
    'D11' = '0.09396'
    'E11' = '0.13%'
    'D12' = '0.001510'
    'E12' = '-0.84%'
    'D13' = '0.04063'
    'E13' = '-12.40%'
    'D14' = '0.0006028'
    'E14' = '0.26%'
    'D15' = '0.006211'
    'E15' = '-0.52%'
    'D16' = '3.513'
    'E16' = '-2.43%'
    'E17' = '-0.55%'
    'D18' = '2.306'
    'E18' = '12.18%'
    'D19' = '0.3149'
    'E19' = '1.21%'
    'D20' = '0.1331'
    'E20' = '-0.25%'
    'D21' = '0.03226'
    'E21' = '1.45%'
    'E22' = '-1.60%'
    'D23' = '3.560'
    'E23' = '-5.51%'
    'D24' = '0.1373'
    'E24' = '-0.10%'
    'D25' = '0.001214'
    'E25' = '-2.41%'
    'D26' = '0.004477'
    'E26' = '-1.66%'
    'D27' = '0.0001179'
    'E27' = '22.77%'
    'E28' = '-27.50%'
    'D40' = '0.03721'
    'E40' = '1.92%'
    'D41' = '0.005956'
    'E41' = '-3.83%'
    'D42' = '0.1056'
    'E42' = '0.38%'
    'D43' = '0.002299'
    'E43' = '-10.72%'
    'D44' = '0.009588'
    'E44' = '20.44%'
    'D45' = '0.00005092'
    'E45' = '-5.58%'
    'E46' = '-0.10%'
    'D47' = '0.1009'
    'E47' = '-30.42%'
    'E48' = '5.25%'
    'D49' = '0.00002099'
    'E49' = '-0.10%'
    'D50' = '0.0001999'
    'E50' = '-0.10%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

